$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.435.10"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +5.64%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.810.86"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +4.28%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.02%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'316.92"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +2.13%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.01%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.5476"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +10.04%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.3862"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +9.16%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'43.10"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +1.48%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.07587"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +4.92%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'1.136"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +7.64%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.05%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'21.16"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +6.40%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'6.227"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +4.83%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = "'1.808.97"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +4.35%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = "'7.318"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +7.13%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'91.18"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +5.76%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'0.00001074"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +4.03%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'0.06479"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +1.52%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.04%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'17.27"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +4.41%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'5.991"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +4.63%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'28.446.53"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +5.42%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'11.32"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +1.48%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'2.126"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +3.85%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = "'20.73"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +4.56%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = "'157.10"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +2.55%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'2.435"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +15.12%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'2.012.29"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +4.04%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'124.22"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +3.15%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'1.169"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +10.99%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'0.1037"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +10.00%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'5.762"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +7.51%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'3.647"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +2.07%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'0.2343"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +17.35%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.02331"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +6.65%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'8.899"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +19.99%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'0.06269"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +6.26%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'11.61"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +5.59%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = "'0.6396"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +6.76%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').Value = "'5.035"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +6.20%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'1.173"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +6.29%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'1.001"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.00%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'1.388"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -2.58%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'13.41"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +4.27%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.6006"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +6.95%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'3.688"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +3.24%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'123.80"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +3.45%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'1.976"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +6.85%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'1.151"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +5.01%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.06936"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +4.25%  "
$ws.Range('E51').Style = 'Normal'
